$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F1").Value = "(* really odd field ☺ *)"
$ws1.Range("F20").Value = "☺ unicode ☺"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("G4").Value = "variance"
$ws2.Range("G5").Value = "std"
$ws2.Range("G9").Value = "coefvar"
$ws2.Range("D17").Select() | Out-Null
$ws1.Activate() | Out-Null
